# "Ban new crude oil plants"
#
# The BBNPPTY ("Ban New Power Plants This Year") sheet has one row per
# power-plant type; each year column (B:AE, years 2021-2050) holds a
# Boolean flag (0/1) saying whether new plants of that type are banned
# starting that year. Row 16 is "crude oil" (see A16). This edit turns
# the ban on for crude oil in every year, and the cells pick up the
# sheet's default/general number format instead of the integer "0"
# format used by the other (still 0/unbanned) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBNPPTY")

$banRow = $ws.Range("B16:AE16")
$banRow.Value = 1
$banRow.Style = "Normal"

# Leave the selection on the cell that was focused when the change was made.
$ws.Range("W7").Select()
